$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.454.10'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.70%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.104.74'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.27%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.89%  '

# Row 6
$ws.Range("E6").Value = '  +0.03%  '

# Row 7
$ws.Range("E7").Value = '  -0.41%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4610'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.61%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.56'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +12.92%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08954'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.46%  '

# Row 11
$ws.Range("E11").Value = '  +1.15%  '

# Row 12
$ws.Range("E12").Value = '  -0.27%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.089.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.52%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.802'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.64%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.935'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.02%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.33%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.004'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.07%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001130'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.53%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06629'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.40%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.99%  '

# Row 22
$ws.Range("E22").Value = '  +0.06%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.515.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.66%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.77%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.365'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.44%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.337.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.45%  '

# Row 27
$ws.Range("E27").Value = '  -0.65%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.566'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.27%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.63'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.06%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.05%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.200'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.00%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1074'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.14%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.693'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.45%  '

# Row 34
$ws.Range("E34").Value = '  +0.44%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.924'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.92%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02572'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.71%  '

# Row 38
$ws.Range("E38").Value = '  +1.87%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.547'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.62%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.90'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.47%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2290'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.92%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6896'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.73%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.248'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.04%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.339'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.80%  '

# Row 45
$ws.Range("E45").Value = '  +0.05%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '14.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.14%  '

# Row 47
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6389'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.27%  '

# Row 48
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.668'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.06%  '

# Row 49
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000359'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +23.89%  '

# Row 51
$ws.Range("E51").Value = '  +1.17%  '
